# Update Karun Nair's per-match batting activity figures (runs/balls/fours)
# for Kings XI Punjab rows, reflecting the latest match added to the form.
# Values are kept as text (as in the source sheet) via the leading
# apostrophe, which forces Excel to store numeric-looking text as text
# instead of converting it to a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (runs/balls/fours)
$ws.Range("C2").Value = "'15"
$ws.Range("D2").Value = "'8"
$ws.Range("E2").Value = "'2"

# Row 3 (runs)
$ws.Range("C3").Value = "'0"

# Row 4 (runs/balls/fours)
$ws.Range("C4").Value = "'1"
$ws.Range("D4").Value = "'3"
$ws.Range("E4").Value = "'0"
